$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2: replace computed date-serial values with literal text dates
# (quote-prefixed so Excel stores them as text, not as date serials)
$ws2.Range("E3").Value = "'12/3/2022"
$ws2.Range("E3").NumberFormat = "mm-dd-yy"

$ws2.Range("E6").Value = "'12/10/2022"
$ws2.Range("E6").NumberFormat = "mm-dd-yy"

$ws2.Range("E9").Value = "'12/17/2022"
$ws2.Range("E9").NumberFormat = "mm-dd-yy"

# Sheet1: same treatment
$ws1.Range("F3").Value = "'10/8/2022"
$ws1.Range("F3").NumberFormat = "mm-dd-yy"

$ws1.Range("F6").Value = "'10/15/2022"
$ws1.Range("F6").NumberFormat = "mm-dd-yy"

$ws1.Range("F8").Value = "'10/29/2022"
$ws1.Range("F8").NumberFormat = "mm-dd-yy"

$ws1.Range("F11").Value = "'11/4/2022"
$ws1.Range("F11").NumberFormat = "mm-dd-yy"

$ws1.Range("F13").Value = "'11/12/2022"
$ws1.Range("F13").NumberFormat = "mm-dd-yy"

$ws1.Range("F15").Value = "'11/26/2022"
$ws1.Range("F15").NumberFormat = "mm-dd-yy"
